# Quality_Assurance_KPI_RH_CH.xlsx - sample template refresh for the
# HP Kayakalp-program download.
#
# The authored "S. No." column (C4:C47) was cleared out so the
# downloaded template starts blank for end users; only the view/selection
# state that Excel persisted on last save changes alongside it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the numbered "S. No." values in column C (rows 4-47) while
# leaving the cell styling (s="4") intact - matches the section-header
# rows (10, 17, 30, 37, 41) that were already blank in that column.
$ws.Range("C4:C47").ClearContents() | Out-Null

# Restore the view state Excel wrote on its last save: scrolled near the
# top of the sheet with D12 selected (previously topLeftCell A31 / B49).
$ws.Range("D12").Select() | Out-Null
